$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header text casing (ID Number -> ID NUMBER)
$ws.Cells.Item(1, 1).Value = "ID NUMBER"

# New data rows: ID Number, Last Name, First Name, Middle Name
$data = @(
    @(20190016812, "ORPILLA", "GBRIEL RENZ", "CABALLEROS"),
    @(20170012416, "PALER", "NATHALIE KATE", $null),
    @(20170011560, "PAPA", "CZIRELLE DOMINIQUE", $null),
    @(20170011644, "LAPUZ", "ZAMANTHA", $null),
    @(20150007960, "DOLENDO", "LEO", $null)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    if ($entry[3] -ne $null) {
        $ws.Cells.Item($row, 4).Value = $entry[3]
    }
    $row++
}

# Column A (the ID Number column) gets a distinct style: Times New Roman, wrap text.
# Build the combined style once on the first data cell, then propagate the already
# resolved style to the remaining cells via a format-only paste so we don't create
# an extra intermediate cell style in the workbook's style table.
$firstIdCell = $ws.Cells.Item(2, 1)
$firstIdCell.Font.Name = "Times New Roman"
$firstIdCell.WrapText = $true

$firstIdCell.Copy()
$restIdRange = $ws.Range("A3:A6")
$restIdRange.PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C7").Select()
